$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 24.01.2022 16:45"

# Row 9 updates: price (B9) and old price (C9) swap values, delta (D9)
# becomes a text "+0.4", and the old-date column (E9) becomes a plain
# text timestamp string. Both D9 and E9 lose their prior numeric styling.
$ws.Range("B9").Value = 37.9
$ws.Range("C9").Value = 37.5

$ws.Range("D9").Value = "'+0.4"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "'2022-01-24 16:45:52"
$ws.Range("E9").Style = "Normal"
